$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.721.10"
$ws.Range("E2").Value = "  +3.71%  "
$ws.Range("D3").Value = "2.423.35"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.01"
$ws.Range("E5").Value = "  +3.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.67"
$ws.Range("E6").Value = "  +6.08%  "
$ws.Range("E7").Value = "  +2.05%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.528"
$ws.Range("E9").Value = "  +9.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.54"
$ws.Range("E10").Value = "  +3.36%  "
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("E12").Value = "  +1.76%  "
$ws.Range("E13").Value = "  -2.58%  "
$ws.Range("E14").Value = "  +2.99%  "
$ws.Range("D15").Value = "2.801.07"
$ws.Range("E15").Value = "  +2.39%  "
$ws.Range("D16").Value = "2.421.77"
$ws.Range("E16").Value = "  +1.67%  "
$ws.Range("E17").Value = "  +4.21%  "
$ws.Range("D18").Value = "44.562.47"
$ws.Range("E18").Value = "  +3.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.33"
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("E20").Value = "  +1.55%  "
$ws.Range("E21").Value = "  +3.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.81"
$ws.Range("E22").Value = "  +0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.47"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("E24").Value = "  +4.47%  "
$ws.Range("E25").Value = "  +2.01%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.22"
$ws.Range("E27").Value = "  +2.68%  "
$ws.Range("E28").Value = "  -3.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.54"
$ws.Range("E29").Value = "  +1.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "33.68"
$ws.Range("E30").Value = "  +4.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "48.55"
$ws.Range("E31").Value = "  +1.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.127"
$ws.Range("E32").Value = "  +19.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.53"
$ws.Range("E33").Value = "  +10.96%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0778"
$ws.Range("E34").Value = "  +8.44%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.17"
$ws.Range("E35").Value = "  +3.21%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("E37").Value = "  +2.56%  "
$ws.Range("E38").Value = "  +3.04%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "121.28"
$ws.Range("E40").Value = "  -2.53%  "
$ws.Range("E41").Value = "  +1.67%  "
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.12"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0291"
$ws.Range("E44").Value = "  +4.25%  "
$ws.Range("D45").Value = "1.946.49"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("E47").Value = "  +8.27%  "
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.69"
$ws.Range("E49").Value = "  +11.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.29"
$ws.Range("E50").Value = "  +5.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "75.15"
$ws.Range("E51").Value = "  +4.61%  "
